$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Append new rows to the "LE" sheet (Ocean Jobs block, rows 10-15)
#    Doing this before creating the "AO" sheet means the new shared
#    string "Ocean Jobs" is registered before the AO-only strings,
#    matching the target shared string order.
# ------------------------------------------------------------------
$wsLE = $wb.Worksheets.Item("LE")

$wsLE.Range("A10").Value = 2013

$wsLE.Range("A11").Value = "Region"
$wsLE.Range("B11").Value = "Ocean Jobs"

$wsLE.Range("A12").Value = "Hawaii"
$wsLE.Range("B12").Value = 13576

$wsLE.Range("A13").Value = "Maui Nui"
$wsLE.Range("B13").Value = 25423

$wsLE.Range("A14").Value = "Oahu"
$wsLE.Range("B14").Value = 59163

$wsLE.Range("A15").Value = "Kauai"
$wsLE.Range("B15").Value = 5264

# New column widths for the added Status/Trend-style columns
$wsLE.Columns.Item(3).ColumnWidth = 14.33
$wsLE.Columns.Item(4).ColumnWidth = 17.83

# ------------------------------------------------------------------
# 2. Add a new worksheet "AO" at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAO = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAO.Name = "AO"

# Fill column A first (region names) - this registers the
# "Kauai & Niʻihau" shared string before "Access"/"Resource" below,
# matching the original authoring order.
$wsAO.Range("A1").Value = "Region"
$wsAO.Range("A2").Value = "Hawaiʻi"
$wsAO.Range("A3").Value = "Maui Nui"
$wsAO.Range("A4").Value = "Oahu"
$wsAO.Range("A5").Value = "Kauai & Niʻihau"

$wsAO.Range("B1").Value = "Access"
$wsAO.Range("C1").Value = "Resource"
$wsAO.Range("D1").Value = "Status"
$wsAO.Range("E1").Value = "Trend"

$wsAO.Range("B2").Value = 0.64
$wsAO.Range("C2").Value = 0.66
$wsAO.Range("D2").Value = 0.65
$wsAO.Range("E2").Value = 0.01

$wsAO.Range("B3").Value = 0.56000000000000005
$wsAO.Range("C3").Value = 0.66
$wsAO.Range("D3").Value = 0.61
$wsAO.Range("E3").Value = -0.01

$wsAO.Range("B4").Value = 0.68
$wsAO.Range("C4").Value = 0.54
$wsAO.Range("D4").Value = 0.61
$wsAO.Range("E4").Value = 0.04

$wsAO.Range("B5").Value = 0.56000000000000005
$wsAO.Range("C5").Value = 0.72
$wsAO.Range("D5").Value = 0.64
$wsAO.Range("E5").Value = 0.02

# ------------------------------------------------------------------
# 3. Fix up selections / active sheet / tab order
# ------------------------------------------------------------------
$wsLE.Activate()
$wsLE.Range("C11").Select()

$wsAO.Activate()
$wsAO.Range("J8").Select()
